# Update "Countries & provincias Spain" COVID table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (shared string index 67 / cell A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 12:52"

# 2) Row 8 -- Ciudad Real
$ws.Range("B8").Value = 5442
$ws.Range("C8").Value = 2532
$ws.Range("D8").Value = 9896
$ws.Range("E8").Value = 585

# 3) Row 11 -- Albacete
$ws.Range("B11").Value = 3506
$ws.Range("C11").Value = 2532
$ws.Range("D11").Value = 9896
$ws.Range("E11").Value = 322

# 4) Row 15 -- Toledo
$ws.Range("B15").Value = 3052
$ws.Range("C15").Value = 2532
$ws.Range("D15").Value = 9896
$ws.Range("E15").Value = 431

# 5) Row 36 -- now Guadalajara (was Castello/Castellon)
$ws.Range("A36").Value = "Guadalajara"
$ws.Range("B36").Value = 1134
$ws.Range("C36").Value = 2532
$ws.Range("D36").Value = 9896
$ws.Range("E36").Value = 153

# 6) Row 37 -- now Castello/Castellon (was Guadalajara)
$ws.Range("A37").Value = "Castello/Castellon"
$ws.Range("B37").Value = 1089
$ws.Range("C37").Value = 246
$ws.Range("D37").Value = 740
$ws.Range("E37").Value = 103

# 7) Row 40 -- now Cuenca (was Avila)
$ws.Range("A40").Value = "Cuenca"
$ws.Range("B40").Value = 920
$ws.Range("C40").Value = 2532
$ws.Range("D40").Value = 9896
$ws.Range("E40").Value = 135

# 8) Row 41 -- now Avila (was Aragon)
$ws.Range("A41").Value = "Avila"
$ws.Range("B41").Value = 917
$ws.Range("C41").Value = 350
$ws.Range("D41").Value = 473
$ws.Range("E41").Value = 94

# 9) Row 42 -- now Aragon (was Cuenca)
$ws.Range("A42").Value = "Aragon"
$ws.Range("B42").Value = 907
$ws.Range("C42").Value = 29
$ws.Range("D42").Value = 838
$ws.Range("E42").Value = 40

# 10) Row 54 -- Melilla, recovered count only
$ws.Range("D54").Value = 81
